# Upload a couple of misc files
# Adds a new row (row 3) to Sheet1 with two more Python2/Python3 migration
# notes, styled like the existing "suppport star expressions" note (B2) but
# with a smaller font size, tweaks the column widths slightly, and moves the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for row 3.
$ws.Range("A3").Value = "anydbm , Queue , thread , StringIO.StringIO , urllib.open"
$ws.Range("B3").Value = "dbm , queue , _thread , io.StringIO , urllib.request.urlopen"

# Reuse the formatting already applied to B2 (WenQuanYi Micro Hei font) for
# the new cells, then shrink the font a bit for the new row.
$ws.Range("B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
$ws.Range("A3:B3").Font.Size = 12

# Minor column width adjustments.
$ws.Range("A1").EntireColumn.ColumnWidth = 56.1666666
$ws.Range("B1").EntireColumn.ColumnWidth = 57.1666666

# Update the active selection.
$ws.Range("A8").Select() | Out-Null
